# "Concept in juiste map gezet & Uren registratie tot vakantie"
#
# Fills in the hours registered for Week 6 (rows 29-31: Woensdag/Donderdag/
# Vrijdag) which were previously blank, and highlights those cells the same
# way the Maandag/Dinsdag rows above them already are (red for Woensdag,
# blue for Donderdag/Vrijdag). Also nudges the saved selection from J21 to
# I21, matching the author's last-touched cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# BGR-encoded colors (OLE_COLOR) matching the workbook's existing
# highlight fills: red (FFFF0000) and blue (FF0070C0).
$red  = 255        # 0x0000FF -> R=FF,G=00,B=00
$blue = 12611584   # 0x C07000 -> R=00,G=70,B=C0

# --- Woensdag (row 29): all zero, highlighted red ---
$ws.Range("B29").Value = 0
$ws.Range("C29").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 0
$ws.Range("C29:G29").Interior.Color = $red

# --- Donderdag (row 30): all one, highlighted blue ---
$ws.Range("B30").Value = 1
$ws.Range("C30").Value = 1
$ws.Range("D30").Value = 1
$ws.Range("E30").Value = 1
$ws.Range("F30").Value = 1
$ws.Range("G30").Value = 1
$ws.Range("C30:G30").Interior.Color = $blue

# --- Vrijdag (row 31): all four, highlighted blue ---
$ws.Range("B31").Value = 4
$ws.Range("C31").Value = 4
$ws.Range("D31").Value = 4
$ws.Range("E31").Value = 4
$ws.Range("F31").Value = 4
$ws.Range("G31").Value = 4
$ws.Range("C31:G31").Interior.Color = $blue

# Totals in row 32 and the summary formulas near the top of the sheet are
# formula-driven and recalc automatically.

# Move the saved selection to I21 (was J21).
$ws.Range("I21").Select()
